# Fill the "duplicate_image_filename" column (E) with "NA" for the
# practice (p1-p4) and generic/unique stimuli rows (rows 2-21), matching
# the other rows that already carry a value for every column.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2:E21").Value = "NA"
